# Rename source_data to data (#157)
#
# The "Data table" sheet is renamed to "Data", and the saved cursor
# selection on that sheet is updated to reflect where the user left off
# after the rename.

$wb = $excel.ActiveWorkbook

# Rename "Data table" sheet to "Data"
$ws = $wb.Worksheets.Item("Data table")
$ws.Name = "Data"

# Restore the active sheet / selection state on the renamed sheet
$ws.Activate()
$ws.Range("H26").Select()
